# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C15").Value = 45233
